# Apply corrections to the Companywise Stock Report.
# A number of stock-adjustment line items had their quantity/value figures
# revised (and, in several cases, two adjacent rows for the same item had
# their Code/Rate2/Qty/Value figures transposed). Subtotal and Grand Total
# rows are updated to reflect the corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BHAWAR SALES CORPORATION block (rows 60-95, subtotal row 96) ---
$ws.Range("F83").Value = 231
$ws.Range("G83").Value = 10796.94
$ws.Range("F92").Value = 42
$ws.Range("G92").Value = 5269.74
$ws.Range("B96").Value = 159752.41

# --- Rows 155/156 (DAB-Real Activ Coconut Water Tetra 1000ml) swapped figures ---
$ws.Range("B155").Value = 64329
$ws.Range("E155").Value = 128.32
$ws.Range("F155").Value = 1
$ws.Range("G155").Value = 120.69

$ws.Range("B156").Value = 57552
$ws.Range("E156").Value = 136.86
$ws.Range("F156").Value = -5
$ws.Range("G156").Value = -603.45

# --- GHP block (rows 176-183, subtotal row 184) ---
$ws.Range("F177").Value = 205
$ws.Range("G177").Value = 13284

$ws.Range("F180").Value = 55
$ws.Range("G180").Value = 4855.95

$ws.Range("B184").Value = 27420.66

# --- Rows 220/221 (HIM-GENTLE BABY SOAP 75G) swapped figures ---
$ws.Range("B220").Value = 48706
$ws.Range("E220").Value = 39.8
$ws.Range("F220").Value = -144
$ws.Range("G220").Value = -4795.2

$ws.Range("B221").Value = 64973
$ws.Range("E221").Value = 35.4
$ws.Range("F221").Value = 0
$ws.Range("G221").Value = 0

# --- Rows 247/248 (HUL-Bru Inst Poly 50g) swapped figures ---
$ws.Range("B247").Value = 63565
$ws.Range("E247").Value = 109.19
$ws.Range("F247").Value = 60
$ws.Range("G247").Value = 6162.6

$ws.Range("B248").Value = 61610
$ws.Range("E248").Value = 122.71
$ws.Range("F248").Value = -58
$ws.Range("G248").Value = -5957.18

# --- HUL block row 253 (HUL-Gal Mv Face Wash 100G) ---
$ws.Range("F253").Value = 70
$ws.Range("G253").Value = 7996.1

# --- Rows 260/261 (HUL-knorr schezwan 200g pch) swapped figures ---
$ws.Range("B260").Value = 55356
$ws.Range("E260").Value = 54.04
$ws.Range("F260").Value = -158
$ws.Range("G260").Value = -7527.12

$ws.Range("B261").Value = 63510
$ws.Range("E261").Value = 50.66
$ws.Range("F261").Value = 112
$ws.Range("G261").Value = 5335.68

# --- Rows 271/272 (Hul-pears pure and gentle 3x125 gm) swapped figures ---
$ws.Range("B271").Value = 60325
$ws.Range("E271").Value = 151.57
$ws.Range("F271").Value = -102
$ws.Range("G271").Value = -12939.72

$ws.Range("B272").Value = 63560
$ws.Range("E272").Value = 134.87
$ws.Range("F272").Value = 1
$ws.Range("G272").Value = 126.86

# --- HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp (row 283) & block subtotal (row 288) ---
$ws.Range("F283").Value = 493
$ws.Range("G283").Value = 84465.69

$ws.Range("B288").Value = 160611.75

# --- Rows 412/413 (CRE-Butter cookies 100gm) swapped figures ---
$ws.Range("B412").Value = 65068
$ws.Range("E412").Value = 13.97
$ws.Range("F412").Value = 63
$ws.Range("G412").Value = 828.45

$ws.Range("B413").Value = 53602
$ws.Range("E413").Value = 15.69
$ws.Range("F413").Value = -231
$ws.Range("G413").Value = -3037.65

# --- Rows 420/421 (CRE-Cremica Golden Bytes Rich Butter 200Gm) swapped figures ---
$ws.Range("B420").Value = 64922
$ws.Range("E420").Value = 20.98
$ws.Range("F420").Value = 67
$ws.Range("G420").Value = 1321.91

$ws.Range("B421").Value = 45706
$ws.Range("E421").Value = 23.58
$ws.Range("F421").Value = -202
$ws.Range("G421").Value = -3985.46

# --- Rows 428/429 (CRE-Cremica Oatmeal Digestive 112.5 Gm) swapped figures ---
$ws.Range("B428").Value = 45709
$ws.Range("E428").Value = 15.69
$ws.Range("F428").Value = -300
$ws.Range("G428").Value = -3945

$ws.Range("B429").Value = 64925
$ws.Range("E429").Value = 13.97
$ws.Range("F429").Value = 111
$ws.Range("G429").Value = 1459.65

# --- Rows 521/522 (Rasna Nagpur Orange (32 Glass)) swapped figures ---
$ws.Range("B521").Value = 60022
$ws.Range("E521").Value = 37.22
$ws.Range("F521").Value = -113
$ws.Range("G521").Value = -3709.79

$ws.Range("B522").Value = 64830
$ws.Range("E522").Value = 34.9
$ws.Range("F522").Value = 101
$ws.Range("G522").Value = 3315.83

# --- VVD block rows 653 & 655, subtotal row 660 ---
$ws.Range("F653").Value = 1190
$ws.Range("G653").Value = 194100.9

$ws.Range("F655").Value = 343
$ws.Range("G655").Value = 97024.41

$ws.Range("B660").Value = 334934.41

# --- Grand Total rows ---
$ws.Range("B679").Value = 3285049.52
$ws.Range("B680").Value = 3285049.52
